# Apply the targeted cell-value edits described in the commit diff.
#
# New SSTH batch number / production data updates:
#   - System Configuration0@0x0!E3 (Radio Name):       BBCB4866          -> BBCB4859
#   - Statistics@0x5!E7            (Production Date):  20191212          -> 20191213
#   - Statistics@0x5!E8            (Batch Number):      98               -> 106
#   - Statistics@0x5!E9            (Serial):            9259266508322    -> 9259266512513
#   - Calibration0@0x8!E3          (Accel X - D):       -99.06767272949219   -> -102.89463806152344
#   - Calibration0@0x8!E9          (Voltage Battery-D): -0.051661375910043716 -> -0.04075469821691513
#
# All of these source cells are stored as plain TEXT (inline strings) even
# though several of the new values look numeric. Excel's COM layer infers a
# Number type for anything that parses as a number unless the cell is
# pre-formatted as Text, so for the numeric-looking values we force the
# "@" (Text) number format before assigning — this keeps the stored value an
# exact-text match (no scientific notation / float rounding) instead of a
# General-formatted number.

$wb = $excel.ActiveWorkbook

# --- System Configuration0@0x0 ---------------------------------------
$wsConfig = $wb.Worksheets.Item("System Configuration0@0x0")
# Alphanumeric value - stays text natively, no format change needed.
$wsConfig.Range("E3").Value = "BBCB4859"

# --- Statistics@0x5 -----------------------------------------------------
$wsStats = $wb.Worksheets.Item("Statistics@0x5")

$prodDate = $wsStats.Range("E7")
$prodDate.NumberFormat = "@"
$prodDate.Value = "20191213"

$batchNum = $wsStats.Range("E8")
$batchNum.NumberFormat = "@"
$batchNum.Value = "106"

$serial = $wsStats.Range("E9")
$serial.NumberFormat = "@"
$serial.Value = "9259266512513"

# --- Calibration0@0x8 -----------------------------------------------------
$wsCal = $wb.Worksheets.Item("Calibration0@0x8")

$accelXD = $wsCal.Range("E3")
$accelXD.NumberFormat = "@"
$accelXD.Value = "-102.89463806152344"

$voltBattD = $wsCal.Range("E9")
$voltBattD.NumberFormat = "@"
$voltBattD.Value = "-0.04075469821691513"
